# Sara-Alert-Format-Invalid-Monitorees.xlsx
# Add three new header columns to the end of the row-1 header row
# (CV1:CX1), matching the new "Race ..." columns that were appended to
# the export format, then leave the selection where the author's last
# cursor move landed.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("CV1").Value = "Race Unknown"
$ws.Range("CW1").Value = "Race Other"
$ws.Range("CX1").Value = "Race Refused to Answer"

# Move the selection/cursor the way the author's session ended up.
$ws.Range("CX6").Select()
